$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Row 58 — 238. Product of Array Except Self
# ---------------------------------------------------------------------------
$ws.Range("A58").Value = "238. Product of Array Except Self"
$ws.Range("B58").Value = "Medium"
$ws.Range("C58").Value = "https://leetcode.com/problems/product-of-array-except-self/"
# Register the hyperlink before the formatting copy below so the copied
# style (which already renders like a hyperlink) wins over the default
# "Hyperlink" style the Add() call would otherwise stamp onto the cell.
$ws.Hyperlinks.Add($ws.Range("C58"), "https://leetcode.com/problems/product-of-array-except-self/")

# Copy formatting (A:F) from row 11, which already carries the exact style
# combination we need (27,8,9,10,11,12); then copy the G-column "not yet
# reviewed" style (44) from row 33's G cell.
$ws.Range("A11:F11").Copy()
$ws.Range("A58:F58").PasteSpecial(-4122)
$ws.Range("G33").Copy()
$ws.Range("G58").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("D58").Value = 44559
$ws.Range("E58").Value = "数学"

$f58 = "类似于分糖果的思路，从左到右遍历一遍，然后再从右往左遍历一遍，记录乘积product，先更新ans，再更新product"
$ws.Range("F58").Value = $f58
$ws.Range("F58").Characters(1, 35).Font.Name = "宋体"
$ws.Range("F58").Characters(36, 7).Font.Name = "Times New Roman"
$ws.Range("F58").Characters(43, 4).Font.Name = "宋体"
$ws.Range("F58").Characters(47, 3).Font.Name = "Times New Roman"
$ws.Range("F58").Characters(50, 4).Font.Name = "宋体"
$ws.Range("F58").Characters(54, 7).Font.Name = "Times New Roman"

$ws.Range("G58").Value = "未复习"

$ws.Rows(58).RowHeight = 42

# ---------------------------------------------------------------------------
# Row 59 — 135. Candy
# ---------------------------------------------------------------------------
$ws.Range("A59").Value = "135. Candy"
$ws.Range("B59").Value = "Hard"
$ws.Range("C59").Value = "https://leetcode.com/problems/candy/"
$ws.Hyperlinks.Add($ws.Range("C59"), "https://leetcode.com/problems/candy/")

# Copy formatting (A:G) from row 57, which already has the exact style
# combination we need (34,18,19,20,21,31,28) without an H cell.
$ws.Range("A57:G57").Copy()
$ws.Range("A59:G59").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("D59").Value = 44442
$ws.Range("E59").Value = "贪心"
$ws.Range("F59").Value = "从左往右遍历使每个元素与相邻右元素满足要求；从右往左遍历使与相邻做元素满足要求"
$ws.Range("G59").Value = 44559

$ws.Rows(59).RowHeight = 28

# ---------------------------------------------------------------------------
# View state — selection moves to F54 (matches the saved selection in the
# target workbook).
# ---------------------------------------------------------------------------
$ws.Range("F54").Select()

Write-Output "done"
